$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.361.07"
$ws.Range("E2").Value = "  -2.35%  "

$ws.Range("D3").Value = "3.535.18"
$ws.Range("E3").Value = "  -4.08%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.35%  "

$ws.Range("D7").Value = "3.526.17"
$ws.Range("E7").Value = "  -3.93%  "

$ws.Range("E8").Value = "  -0.90%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("E10").Value = "  -4.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.588"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.10%  "

$ws.Range("E14").Value = "  -4.17%  "

$ws.Range("D15").Value = "4.092.87"
$ws.Range("E15").Value = "  -4.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "629.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.12%  "

$ws.Range("D18").Value = "3.535.07"
$ws.Range("E18").Value = "  -4.34%  "

$ws.Range("D19").Value = "69.325.03"
$ws.Range("E19").Value = "  -2.65%  "

$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.892"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.14%  "

$ws.Range("E24").Value = "  -7.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.85%  "

$ws.Range("E26").Value = "  -3.84%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("E28").Value = "  -6.60%  "

$ws.Range("E29").Value = "  -8.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.10%  "

$ws.Range("E31").Value = "  -7.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.90%  "

$ws.Range("E33").Value = "  -6.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "634.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -13.28%  "

$ws.Range("E38").Value = "  -4.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0456"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.19%  "

$ws.Range("E42").Value = "  -5.05%  "

$ws.Range("D43").Value = "3.394.36"
$ws.Range("E43").Value = "  -6.25%  "

$ws.Range("E44").Value = "  -5.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "33.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.26%  "

$ws.Range("D46").Value = "0.0₃0701"
$ws.Range("E46").Value = "  -8.45%  "

$ws.Range("E47").Value = "  -6.73%  "

$ws.Range("E48").Value = "  -3.90%  "

$ws.Range("E49").Value = "  -2.23%  "

$ws.Range("E50").Value = "  +14.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.87%  "
